$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166208744049072
$ws.Range("B1").Value = 2.427464485168457
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.372320413589478
$ws.Range("E1").Value = 1.234760999679565
